$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K3").Value = 2023
$ws.Range("K4").Value = 534
$ws.Range("K5").Value = 186
$ws.Range("K6").Value = 348
